$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "work and ses data" - two instance counts changed on the roll-value table
$ws.Range("B8").Value = 6
$ws.Range("B17").Value = 10

# Selection moved (mirrors the saved cursor position in the source file)
$ws.Range("B9").Select()
